# Updated symbol list on Mon Feb 13 19:21:14 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (D) and "Volume(1h)" (E) columns for each crypto
# symbol row with the latest scraped values. Values are written with a
# leading apostrophe so Excel stores them as literal text (matching the
# original inlineStr cells) rather than re-parsing them as numbers, which
# would otherwise strip significant trailing zeros (e.g. "1.520" -> 1.52)
# or flip tiny magnitudes into scientific notation (e.g. "0.00000000751").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'288.89"
$ws.Range("E2").Value = "'-9.60%"
$ws.Range("D3").Value = "'40.27"
$ws.Range("E3").Value = "'-2.63%"
$ws.Range("D4").Value = "'5.045"
$ws.Range("E4").Value = "'-4.05%"
$ws.Range("D5").Value = "'0.07295"
$ws.Range("E5").Value = "'-5.70%"
$ws.Range("D6").Value = "'4.286"
$ws.Range("E6").Value = "'-1.29%"
$ws.Range("D7").Value = "'1.520"
$ws.Range("E7").Value = "'-12.84%"
$ws.Range("D8").Value = "'0.9187"
$ws.Range("E8").Value = "'-2.66%"
$ws.Range("D9").Value = "'0.1189"
$ws.Range("E9").Value = "'-4.38%"
$ws.Range("D10").Value = "'0.1729"
$ws.Range("E10").Value = "'-7.56%"
$ws.Range("D11").Value = "'0.08635"
$ws.Range("E11").Value = "'-6.58%"
$ws.Range("D12").Value = "'0.04172"
$ws.Range("E12").Value = "'-3.29%"
$ws.Range("D13").Value = "'0.1054"
$ws.Range("E13").Value = "'0.26%"
$ws.Range("D14").Value = "'0.001263"
$ws.Range("D15").Value = "'0.005912"
$ws.Range("E15").Value = "'-0.50%"
$ws.Range("D16").Value = "'3.397"
$ws.Range("E16").Value = "'1.66%"
$ws.Range("E17").Value = "'-1.15%"
$ws.Range("D18").Value = "'0.3291"
$ws.Range("E18").Value = "'-1.24%"
$ws.Range("D19").Value = "'7.840"
$ws.Range("E19").Value = "'1.34%"
$ws.Range("D20").Value = "'0.1351"
$ws.Range("E20").Value = "'-0.25%"
$ws.Range("D21").Value = "'0.2883"
$ws.Range("E21").Value = "'2.00%"
$ws.Range("D22").Value = "'0.03861"
$ws.Range("E22").Value = "'-4.44%"
$ws.Range("D23").Value = "'0.001268"
$ws.Range("E23").Value = "'-0.12%"
$ws.Range("D24").Value = "'0.003837"
$ws.Range("E24").Value = "'-6.71%"
$ws.Range("E25").Value = "'0.63%"
$ws.Range("D26").Value = "'0.0003724"
$ws.Range("D38").Value = "'0.02317"
$ws.Range("E38").Value = "'-9.18%"
$ws.Range("D39").Value = "'0.04973"
$ws.Range("E39").Value = "'-6.89%"
$ws.Range("D40").Value = "'0.006353"
$ws.Range("E40").Value = "'218.69%"
$ws.Range("D41").Value = "'0.007678"
$ws.Range("E41").Value = "'-1.29%"
$ws.Range("E42").Value = "'-3.34%"
$ws.Range("D43").Value = "'0.007348"
$ws.Range("E43").Value = "'0.73%"
$ws.Range("D44").Value = "'0.007068"
$ws.Range("E44").Value = "'-14.64%"
$ws.Range("D45").Value = "'0.3122"
$ws.Range("E45").Value = "'-1.71%"
$ws.Range("D46").Value = "'0.00006424"
$ws.Range("E46").Value = "'-3.82%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.10%"
$ws.Range("D48").Value = "'0.01075"
$ws.Range("E48").Value = "'-94.66%"
$ws.Range("E49").Value = "'-0.14%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.10%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.10%"
